# "Atualizacoes dados 16/07 23h"
# Add a new "dt_insertion" column (G) to the tournament sheet: a header
# cell matching the existing header style, and a timestamp value for the
# single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G1: header cell, formatted like the other header cells (bold font,
#     thin border, centered/top aligned) -------------------------------
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value = "dt_insertion"

# --- G2: insertion timestamp for the existing data row -----------------
$ws.Range("G2").Value = 45489.94324074074
$ws.Range("G2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("G2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
